$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2018-12-31 00:00:00"
$ws.Range("O2").Value = 45724972.05
$ws.Range("P2").Value = 204779402.99
$ws.Range("Q2").Value = 150709675.39
$ws.Range("R2").Value = 45.7001254974
$ws.Range("S2").Value = 115922904.19
$ws.Range("T2").Value = 115922904.19
$ws.Range("U2").Value = 52.6115040239
$ws.Range("V2").Value = 4995910.32
$ws.Range("W2").Value = 13495810.32
$ws.Range("X2").Value = 245758.04
$ws.Range("Y2").Value = 51878341.98
$ws.Range("Z2").Value = 52471127.33
$ws.Range("AA2").Value = 6746155.28
$ws.Range("AG2").Value = 1402040.97
$ws.Range("AP2").Value = 49.5042043426
$ws.Range("AQ2").Value = 46.354084362197
$ws.Range("AR2").Value = 51.74842211189
$ws.Range("AS2").Value = 44941740.83
$ws.Range("AT2").Value = 63.711150925553
